$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each row below gives the (row, col) address of a table cell together
# with the text currently in it and the text it must contain afterwards.
# Addressing cells directly (rather than a document-wide Find/Replace) is
# required because some source values are duplicated elsewhere in the
# table (e.g. "916÷7=130, 6" occurs twice, changing to two different
# results), so a plain text search-and-replace could not tell them apart.
$edits = @(
    @{ Row = 1; Col = 1; Old = "408÷8=51, 0"; New = "596÷8=74, 4" }
    @{ Row = 1; Col = 2; Old = "324÷2=162, 0"; New = "339÷7=48, 3" }
    @{ Row = 1; Col = 3; Old = "634÷7=90, 4"; New = "901÷5=180, 1" }
    @{ Row = 1; Col = 4; Old = "716÷7=102, 2"; New = "226÷4=56, 2" }
    @{ Row = 1; Col = 5; Old = "513÷4=128, 1"; New = "462÷2=231, 0" }
    @{ Row = 5; Col = 1; Old = "181÷4=45, 1"; New = "634÷3=211, 1" }
    @{ Row = 5; Col = 2; Old = "615÷9=68, 3"; New = "125÷2=62, 1" }
    @{ Row = 5; Col = 3; Old = "708÷4=177, 0"; New = "266÷6=44, 2" }
    @{ Row = 5; Col = 4; Old = "480÷7=68, 4"; New = "255÷7=36, 3" }
    @{ Row = 5; Col = 5; Old = "512÷7=73, 1"; New = "384÷3=128, 0" }
    @{ Row = 9; Col = 1; Old = "883÷7=126, 1"; New = "546÷3=182, 0" }
    @{ Row = 9; Col = 2; Old = "503÷4=125, 3"; New = "941÷6=156, 5" }
    @{ Row = 9; Col = 3; Old = "281÷6=46, 5"; New = "275÷6=45, 5" }
    @{ Row = 9; Col = 4; Old = "302÷5=60, 2"; New = "963÷8=120, 3" }
    @{ Row = 9; Col = 5; Old = "473÷6=78, 5"; New = "350÷4=87, 2" }
    @{ Row = 13; Col = 1; Old = "345÷7=49, 2"; New = "761÷7=108, 5" }
    @{ Row = 13; Col = 2; Old = "916÷7=130, 6"; New = "311÷7=44, 3" }
    @{ Row = 13; Col = 3; Old = "376÷7=53, 5"; New = "676÷8=84, 4" }
    @{ Row = 13; Col = 4; Old = "242÷2=121, 0"; New = "262÷9=29, 1" }
    @{ Row = 13; Col = 5; Old = "947÷5=189, 2"; New = "526÷6=87, 4" }
    @{ Row = 17; Col = 1; Old = "732÷7=104, 4"; New = "869÷9=96, 5" }
    @{ Row = 17; Col = 2; Old = "300÷7=42, 6"; New = "233÷6=38, 5" }
    @{ Row = 17; Col = 3; Old = "916÷7=130, 6"; New = "263÷5=52, 3" }
    @{ Row = 17; Col = 4; Old = "524÷7=74, 6"; New = "839÷3=279, 2" }
    @{ Row = 17; Col = 5; Old = "848÷7=121, 1"; New = "310÷3=103, 1" }
)

foreach ($edit in $edits) {
    $cell = $t.Cell($edit.Row, $edit.Col)
    $range = $cell.Range
    # The cell Range.Text includes a trailing cell-mark character, so
    # compare using StartsWith rather than an exact match.
    if (-not $range.Text.StartsWith($edit.Old)) {
        throw "Unexpected text in cell ($($edit.Row), $($edit.Col)): [$($range.Text)] expected [$($edit.Old)]"
    }
    $range.Text = $edit.New
}

Write-Host "Updated $($edits.Count) cells"
